$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 17) mirroring the format of the existing rows.
$row = 17

$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 42622.888495370367

$ws.Cells.Item($row, 2).Value = -10
$ws.Cells.Item($row, 3).Value = 58
$ws.Cells.Item($row, 4).Value = 39
$ws.Cells.Item($row, 5).Value = 58
$ws.Cells.Item($row, 6).Value = 27
$ws.Cells.Item($row, 7).Value = 10384
$ws.Cells.Item($row, 8).Value = 18782
$ws.Cells.Item($row, 9).Value = 2082
$ws.Cells.Item($row, 10).Value = 293
$ws.Cells.Item($row, 11).Value = 198
$ws.Cells.Item($row, 12).Value = 35
$ws.Cells.Item($row, 13).Value = 13
$ws.Cells.Item($row, 14).Value = "Bag"
